$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit cyclically rotates a subset of fields among rows 15, 16 and 17:
#   row15 <- row16's old values
#   row16 <- row17's old values
#   row17 <- row15's old values
# across columns A, P, Q, R, Z, AB, AF, AW, AX.

$rows = @(15, 16, 17)
$cols = @("A", "P", "Q", "R", "Z", "AB", "AW", "AX")

# Capture the current ("before") values for every affected cell using
# Value2 (the plain `.Value` getter isn't usable through this bridge).
$old = @{}
foreach ($r in $rows) {
    foreach ($col in $cols) {
        $addr = "$col$r"
        $old[$addr] = $ws.Range($addr).Value2
    }
}

# AF is handled separately because a blank cell and a wholly absent cell
# both read back as empty through Value2 - tell them apart first.
$afVal = @{}
$afPresent = @{}
foreach ($r in $rows) {
    $v = $ws.Range("AF$r").Value2
    $afVal[$r] = $v
    $afPresent[$r] = ($v -ne $null)
}

# Source row for each destination row (cyclic rotation).
$srcRow = @{ 15 = 16; 16 = 17; 17 = 15 }

foreach ($r in $rows) {
    $s = $srcRow[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $old["$col$s"]
    }
}

foreach ($r in $rows) {
    $s = $srcRow[$r]
    if ($afPresent[$s]) {
        if ($afVal[$s] -eq $null -or $afVal[$s] -eq "") {
            # Source cell existed but was blank. Leave destinations that
            # are already present-and-blank untouched (no real change),
            # and only force-create a cell where one didn't exist before -
            # plain Value="" would just drop it instead of leaving a
            # present-but-blank cell.
            if (-not $afPresent[$r]) {
                $ws.Range("AF$r").Value = ""
                $ws.Range("AF$r").Style = "Normal"
            }
        } else {
            $ws.Range("AF$r").Value = $afVal[$s]
        }
    } else {
        # Source cell was wholly absent - destination must end up absent.
        if ($afPresent[$r]) {
            $ws.Range("AF$r").Value = ""
        }
    }
}
